$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the base value in B1
$ws.Range("B1").Value = 10

# Update the formula in C2 (was "=$B$1*2", now "=$B$1")
$ws.Range("C2").Formula = "=`$B`$1"

# Move the active selection to C3
$ws.Range("C3").Select()
